# 4_analyze/app.R with downloading of plots and tmas
# Adds a new row (row 5) describing the 4_analyze step to the workflow
# summary table, widens column A to fit the longer "Input" text, and
# grows the row heights of the wrapped-text rows to match the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new row of data -------------------------------------------
# Written in this column order (B, D, E, C, F, A) so that new shared
# strings are interned in the same order as the source workbook.
$ws.Range("B5").Value = "4_analyze/app.R"
$ws.Range("D5").Value = "Robust dRFU and curve fitting analysis by multiple possible fitting models"
$ws.Range("E5").Value = "Tmas for input data, and sigmoid fits "
$ws.Range("C5").Value = "analysis.R"
$ws.Range("F5").Value = "values`$tm_table_dRFU, values`$df_models, values`$df_tm_models, values`$df_BIC_display "
$ws.Range("A5").Value = "formatted data, with or without layout (values`$df)"

# --- Widen column A so the longer "Input" entries read comfortably -----
$ws.Columns.Item(1).ColumnWidth = 43.33

# --- Grow the wrapped-text row heights for the now-taller cell content -
$ws.Rows.Item(3).RowHeight = 48
$ws.Rows.Item(4).RowHeight = 48
$ws.Rows.Item(5).RowHeight = 64

# --- Update the active selection ----------------------------------------
$ws.Range("A10:XFD15").Select() | Out-Null
